$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 21: new data row ---
$ws.Range("A21").Value = 45258.681944444441
$ws.Range("A21").NumberFormat = "m/d/yy h:mm"
$ws.Range("B21").Value = 1776128
$ws.Range("C21").Value = 693760
$ws.Range("D21").Value = 392704
$ws.Range("E21").Formula = "=SUM(B21:D21)"
$ws.Range("F21").Value = 170513693
$ws.Range("G21").Formula = "=1-(E21/E2)"
$ws.Range("H21").Formula = "=1-(F21/F2)"
$ws.Range("I21").Value = "Remove software rendering menu and cvars"

# --- Row 22: new data row ---
$ws.Range("A22").Value = 45258.779861111114
$ws.Range("A22").NumberFormat = "m/d/yy h:mm"
$ws.Range("B22").Value = 1776128
$ws.Range("C22").Value = 693760
$ws.Range("D22").Value = 392704
$ws.Range("E22").Formula = "=SUM(B22:D22)"
$ws.Range("F22").Value = 151740398
$ws.Range("G22").Formula = "=1-(E22/E2)"
$ws.Range("H22").Formula = "=1-(F22/F2)"
$ws.Range("I22").Value = "Remove duplicate TGA textures"

# --- Selection moves to C8 ---
$ws.Range("C8").Select()

$wb.Save()
